$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the new data area (column F helper values + summary rows) ---
$ws.Range("F1").Value = 1
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 1

$ws.Range("A3").Value = "cell with formula  condition"
$ws.Range("B3").Formula = "=SUBTOTAL(109,F1:F5)"

$ws.Range("A4").Value = "cell with formula value and > condition"
$ws.Range("B4").Value = 5

# --- Widen column A so the new labels are fully visible ---
$ws.Columns.Item(1).ColumnWidth = 31.33

# --- New conditional formatting rules ---
# B4: highlight when value is less than the subtotal in B3
$fc1 = $ws.Range("B4").FormatConditions.Add(1, 6, "=`$B`$3")
$fc1.Font.Color = 393372
$fc1.Interior.Color = 13551615
$fc1.Priority = 2

# B3: highlight when the subtotal formula result is greater than 5
$fc2 = $ws.Range("B3").FormatConditions.Add(2, 0, "=`$B`$3>5")
$fc2.Font.Color = 393372
$fc2.Interior.Color = 13551615
$fc2.Priority = 1

# Push the pre-existing rules down in precedence (same as Excel does when new
# rules are inserted above older ones).
$a1Rule = $ws.Range("A1").FormatConditions.Item(1)
$a1Rule.Priority = 4
$b1Rule = $ws.Range("B1").FormatConditions.Item(1)
$b1Rule.Priority = 5

# --- Leave the active selection on B4, matching the saved view state ---
$null = $ws.Range("B4").Select()
